# Astro Legends review edit:
#  1. Remove the "Meta description: ..." paragraph that follows the title.
#  2. Replace the final "Create a cartoon-style feature image..." paragraph with
#     two paragraphs: a bold title line followed by the (former) meta-description
#     text, now in italics.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: delete the "Meta description" paragraph (the 2nd paragraph, right
# after the H1 title).
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
if ($metaPara.Range.Text -notmatch "^Meta description") {
    throw "Unexpected document structure: paragraph 2 is not the Meta description paragraph."
}
[void]$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# Part 2: turn the last paragraph (currently the "Create a cartoon-style..."
# image-prompt paragraph) into two paragraphs:
#   - a new bold paragraph with the page title text
#   - the existing paragraph, with its text replaced (keeping italics)
# ---------------------------------------------------------------------------

# Insert a new empty paragraph right before the last paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
[void]$lastPara.Range.InsertParagraphBefore()

# Fill that new (now second-to-last) paragraph with the bold title text via
# an OOXML fragment, so the run formatting matches exactly (a leading empty
# run followed by a bold run), instead of inheriting formatting from Font.*.
$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:b/></w:rPr>
              <w:t>Play Astro Legends: Lyra and Erion Slot Free | Review 2021</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$newParaRange.InsertXML($titleXml)

# Replace the text of the now-last paragraph (the old "Create a cartoon-style
# ..." image prompt) with the former meta-description sentence, keeping the
# paragraph's existing leading empty run and its italic formatting.
$imgPara = $d.Paragraphs($d.Paragraphs.Count)
if ($imgPara.Range.Text -notmatch "^Create a cartoon-style") {
    throw "Unexpected document structure: last paragraph is not the image-prompt paragraph."
}
$imgParaRange = $d.Range($imgPara.Range.Start, $imgPara.Range.End)
$descXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr><w:i/></w:rPr>
              <w:t>Discover the unique mechanics and retro sci-fi theme of Astro Legends in this 2021 review. Play for free and learn about the Lyra Spirit Bonus and Wild Multipliers.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
[void]$imgParaRange.InsertXML($descXml)
